# Germany Landesliga - database update (17-03-2024)
# 1) Rows 14/15 had their match data in the wrong order - swap the two
#    matches' details back into the correct rows.
# 2) Append four new finished/scheduled matches as rows 66-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows 14 and 15 (swap match data) ---
# Row 14
$ws.Range("B14").Value = 7089910
$ws.Range("F14").Value = "ASV Mettmann"
$ws.Range("G14").Value = "TuRU Dsseldorf"
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = 3.25
$ws.Range("L14").Value = 4
$ws.Range("M14").Value = 1.8
$ws.Range("N14").Value = 3.25
$ws.Range("O14").Value = 4
$ws.Range("P14").Value = 1.8
$ws.Range("Q14").Value = 0.5
$ws.Range("R14").Value = 1.975
$ws.Range("S14").Value = 1.825
$ws.Range("T14").Value = 3.25
$ws.Range("U14").Value = 1.85
$ws.Range("V14").Value = 1.95
$ws.Range("W14").Value = 2.25
$ws.Range("Z14").Value = 0.9750000000000001
$ws.Range("AB14").Value = -0.5
$ws.Range("AC14").Value = 0.475

# Row 15
$ws.Range("B15").Value = 7089911
$ws.Range("F15").Value = "Spvgg Steele 0309"
$ws.Range("G15").Value = "VfB Frohnhausen"
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 2.25
$ws.Range("L15").Value = 3.75
$ws.Range("M15").Value = 2.5
$ws.Range("N15").Value = 2.25
$ws.Range("O15").Value = 3.75
$ws.Range("P15").Value = 2.5
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 1.8
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 3.5
$ws.Range("U15").Value = 1.8
$ws.Range("V15").Value = 2
$ws.Range("W15").Value = 1.25
$ws.Range("Z15").Value = 0.8
$ws.Range("AB15").Value = 0.8
$ws.Range("AC15").Value = -1

# --- Append new match rows 66-69 ---
# Row 66
$ws.Range("A66").Value = 64
$ws.Range("B66").Value = 7951581
$ws.Range("C66").Value = "Germany Landesliga"
$ws.Range("D66").Value = "Germany Landesliga"
$ws.Range("E66").Value = 45363.66666666666
$ws.Range("F66").Value = "FC Concordia 03"
$ws.Range("G66").Value = "Germania Schneiche"
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 3
$ws.Range("J66").Value = "A"
$ws.Range("K66").Value = 2.9
$ws.Range("L66").Value = 4
$ws.Range("M66").Value = 1.909
$ws.Range("N66").Value = 2.875
$ws.Range("O66").Value = 4
$ws.Range("P66").Value = 1.909
$ws.Range("Q66").Value = 0.5
$ws.Range("R66").Value = 1.825
$ws.Range("S66").Value = 1.975
$ws.Range("T66").Value = 3.75
$ws.Range("U66").Value = 1.925
$ws.Range("V66").Value = 1.875
$ws.Range("W66").Value = -1
$ws.Range("X66").Value = -1
$ws.Range("Y66").Value = 0.909
$ws.Range("Z66").Value = -1
$ws.Range("AA66").Value = 0.9750000000000001
$ws.Range("AB66").Value = -1
$ws.Range("AC66").Value = 0.875

# Row 67
$ws.Range("A67").Value = 65
$ws.Range("B67").Value = 7950775
$ws.Range("C67").Value = "Germany Landesliga"
$ws.Range("D67").Value = "Germany Landesliga"
$ws.Range("E67").Value = 45364.66666666666
$ws.Range("F67").Value = "TV Dinklage"
$ws.Range("G67").Value = "Viktoria Gesmold"
$ws.Range("H67").Value = 4
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = "H"
$ws.Range("K67").Value = 2.15
$ws.Range("L67").Value = 3.4
$ws.Range("M67").Value = 2.8
$ws.Range("N67").Value = 2.25
$ws.Range("O67").Value = 3.5
$ws.Range("P67").Value = 2.625
$ws.Range("Q67").Value = 0
$ws.Range("R67").Value = 1.725
$ws.Range("S67").Value = 2.075
$ws.Range("T67").Value = 3
$ws.Range("U67").Value = 1.8
$ws.Range("V67").Value = 2
$ws.Range("W67").Value = 1.25
$ws.Range("X67").Value = -1
$ws.Range("Y67").Value = -1
$ws.Range("Z67").Value = 0.7250000000000001
$ws.Range("AA67").Value = -1
$ws.Range("AB67").Value = 0.8
$ws.Range("AC67").Value = -1

# Row 68
$ws.Range("A68").Value = 66
$ws.Range("B68").Value = 7969732
$ws.Range("C68").Value = "Germany Landesliga"
$ws.Range("D68").Value = "Germany Landesliga"
$ws.Range("E68").Value = 45368.45833333334
$ws.Range("F68").Value = "FC Monheim"
$ws.Range("G68").Value = "TuRU Dsseldorf"
$ws.Range("K68").Value = 1.333
$ws.Range("L68").Value = 5
$ws.Range("M68").Value = 6
$ws.Range("N68").Value = 1.285
$ws.Range("O68").Value = 5
$ws.Range("P68").Value = 7.5
$ws.Range("Q68").Value = -1.5
$ws.Range("R68").Value = 1.775
$ws.Range("S68").Value = 2.025
$ws.Range("T68").Value = 3
$ws.Range("U68").Value = 1.775
$ws.Range("V68").Value = 2.025
$ws.Range("W68").Value = 0
$ws.Range("X68").Value = 0
$ws.Range("Y68").Value = 0
$ws.Range("Z68").Value = 0
$ws.Range("AA68").Value = 0

# Row 69
$ws.Range("A69").Value = 67
$ws.Range("B69").Value = 7969733
$ws.Range("C69").Value = "Germany Landesliga"
$ws.Range("D69").Value = "Germany Landesliga"
$ws.Range("E69").Value = 45368.47916666666
$ws.Range("F69").Value = "SC Dsseldorf West"
$ws.Range("G69").Value = "DV Solingen"
$ws.Range("K69").Value = 3.9
$ws.Range("L69").Value = 4.5
$ws.Range("M69").Value = 1.571
$ws.Range("N69").Value = 3.8
$ws.Range("O69").Value = 4.5
$ws.Range("P69").Value = 1.571
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = 1.875
$ws.Range("S69").Value = 1.925
$ws.Range("T69").Value = 3.5
$ws.Range("U69").Value = 1.775
$ws.Range("V69").Value = 2.025
$ws.Range("W69").Value = 0
$ws.Range("X69").Value = 0
$ws.Range("Y69").Value = 0
$ws.Range("Z69").Value = 0
$ws.Range("AA69").Value = 0

# --- Apply formatting to new rows (copy from row 65) ---
$ws.Range("A65").Copy() | Out-Null
$ws.Range("A66:A69").PasteSpecial(-4122) | Out-Null
$ws.Range("E65").Copy() | Out-Null
$ws.Range("E66:E69").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
